$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values remain stored as text (matching source formatting)
$textCells = @("D5","D6","D9","D11","D13","D16","D18","D19","D20","D21","D24","D25","D27","D28","D34","D36","D37","D40","D42","D44","D47","D49","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values from the latest crypto data refresh
$ws.Range("D2").Value = '63.709.62'
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("D3").Value = '2.616.88'
$ws.Range("E3").Value = '  -0.16%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '593.95'
$ws.Range("E5").Value = '  -0.17%  '
$ws.Range("D6").Value = '150.88'
$ws.Range("E6").Value = '  +0.73%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -0.21%  '
$ws.Range("D9").Value = '0.113'
$ws.Range("E9").Value = '  +3.96%  '
$ws.Range("E10").Value = '  +3.56%  '
$ws.Range("D11").Value = '0.394'
$ws.Range("E11").Value = '  +3.12%  '
$ws.Range("E12").Value = '  +0.93%  '
$ws.Range("D13").Value = '27.91'
$ws.Range("E13").Value = '  +1.52%  '
$ws.Range("D14").Value = '3.087.25'
$ws.Range("E14").Value = '  -0.13%  '
$ws.Range("D15").Value = '63.562.76'
$ws.Range("E15").Value = '  -0.25%  '
$ws.Range("D16").Value = '0.0000165'
$ws.Range("E16").Value = '  +11.09%  '
$ws.Range("D17").Value = '2.637.46'
$ws.Range("E17").Value = '  +0.00%  '
$ws.Range("D18").Value = '12.24'
$ws.Range("E18").Value = '  +0.82%  '
$ws.Range("D19").Value = '4.80'
$ws.Range("E19").Value = '  +3.94%  '
$ws.Range("D20").Value = '348.66'
$ws.Range("E20").Value = '  +0.08%  '
$ws.Range("D21").Value = '7.00'
$ws.Range("E21").Value = '  +1.84%  '
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("E23").Value = '  +1.64%  '
$ws.Range("D24").Value = '1.68'
$ws.Range("E24").Value = '  -2.03%  '
$ws.Range("D25").Value = '9.24'
$ws.Range("E25").Value = '  +0.66%  '
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("D27").Value = '8.44'
$ws.Range("E27").Value = '  +4.43%  '
$ws.Range("D28").Value = '547.88'
$ws.Range("E28").Value = '  +0.81%  '
$ws.Range("E29").Value = '  -0.34%  '
$ws.Range("E30").Value = '  -0.15%  '
$ws.Range("E31").Value = '  +1.53%  '
$ws.Range("D32").Value = '0.0₃0894'
$ws.Range("E32").Value = '  +5.61%  '
$ws.Range("E33").Value = '  +2.41%  '
$ws.Range("D34").Value = '5.45'
$ws.Range("E34").Value = '  +4.39%  '
$ws.Range("E35").Value = '  +1.38%  '
$ws.Range("B36").Value = 'PolygonEcosystemToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D36").Value = '0.418'
$ws.Range("E36").Value = '  +3.15%  '
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").Value = '164.43'
$ws.Range("E37").Value = '  -2.07%  '
$ws.Range("E38").Value = '  +1.70%  '
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("D40").Value = '19.73'
$ws.Range("E40").Value = '  +1.92%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("D42").Value = '166.96'
$ws.Range("E42").Value = '  -0.80%  '
$ws.Range("E43").Value = '  +4.69%  '
$ws.Range("D44").Value = '23.47'
$ws.Range("E44").Value = '  +9.91%  '
$ws.Range("E45").Value = '  -0.53%  '
$ws.Range("E46").Value = '  +8.60%  '
$ws.Range("D47").Value = '0.636'
$ws.Range("E47").Value = '  +1.35%  '
$ws.Range("E48").Value = '  +3.02%  '
$ws.Range("D49").Value = '0.0965'
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("D50").Value = '19.32'
$ws.Range("E50").Value = '  +1.04%  '
$ws.Range("D51").Value = '0.0₆0231'
$ws.Range("E51").Value = '  +18.06%  '
